$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Config")
$ws2 = $wb.Worksheets.Item("EmailSettings")

# --- EmailSettings sheet content updates ---

# Drop the old hyperlink that lived on A2 (mailto:sibesh@zanui.com)
foreach ($hl in $ws2.Hyperlinks) {
    $hl.Delete()
}

# D2: new reply-to address, now carries the hyperlink (and the Hyperlink style
# that used to live on A2).
$ws2.Range("D2").Value = "noreply@isalbi.com"
$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:noreply@isalbi.com")
$ws2.Range("D2").Style = "Hyperlink"

# A2: new sender address, no longer a hyperlink - styled with a small grey
# Roboto font instead. Build the font on a scratch cell first and copy the
# format across so A2 doesn't retain any trace of the old Hyperlink style.
$scratch = $ws2.Range("Z1")
$scratch.Font.Name = "Roboto"
$scratch.Font.Size = 8
$scratch.Font.Color = 6841183
$scratch.Copy()

$ws2.Range("A2").Value = "lalita.kashyaponestop@gmail.com"
$ws2.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$scratch.Clear()

# Widen column A to fit the new (longer) address
$ws2.Columns.Item(1).ColumnWidth = 21.9440104166667

# Portrait page orientation for EmailSettings
$ws2.PageSetup.Orientation = 1

# --- Selection / active sheet bookkeeping ---

$ws1.Activate()
$ws1.Range("B2").Select()

$ws2.Activate()
$ws2.Range("D3").Select()
